$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "38.118.21"
$ws.Range("E2").Value = "  +2.99%  "

Set-TextValue "D3" "2.061.21"
$ws.Range("E3").Value = "  +2.64%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("E5").Value = "  +1.93%  "

$ws.Range("E6").Value = "  +1.61%  "

Set-TextValue "D7" "58.34"
$ws.Range("E7").Value = "  +6.91%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +2.48%  "

$ws.Range("E10").Value = "  +2.70%  "

$ws.Range("E11").Value = "  -1.10%  "

Set-TextValue "D12" "2.366.07"
$ws.Range("E12").Value = "  +2.71%  "

$ws.Range("E13").Value = "  +3.58%  "

Set-TextValue "D14" "20.70"
$ws.Range("E14").Value = "  +2.54%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D15" "0.755"
$ws.Range("E15").Value = "  +2.20%  "

$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D16" "5.30"
$ws.Range("E16").Value = "  +3.80%  "

Set-TextValue "D17" "2.050.17"
$ws.Range("E17").Value = "  +2.14%  "

Set-TextValue "D18" "38.048.03"
$ws.Range("E18").Value = "  +2.86%  "

Set-TextValue "D19" "6.19"
$ws.Range("E19").Value = "  +1.60%  "

Set-TextValue "D20" "69.75"
$ws.Range("E20").Value = "  +1.60%  "

$ws.Range("E21").Value = "  +1.97%  "

Set-TextValue "D22" "224.69"
$ws.Range("E22").Value = "  +0.78%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("E25").Value = "  +3.24%  "

$ws.Range("E26").Value = "  +2.19%  "

Set-TextValue "D27" "165.90"
$ws.Range("E27").Value = "  +0.30%  "

Set-TextValue "D28" "0.133"
$ws.Range("E28").Value = "  +7.11%  "

Set-TextValue "D29" "19.07"
$ws.Range("E29").Value = "  +2.45%  "

$ws.Range("E30").Value = "  +2.53%  "

$ws.Range("E31").Value = "  +1.87%  "

$ws.Range("E32").Value = "  +1.42%  "

$ws.Range("E33").Value = "  +4.70%  "

Set-TextValue "D34" "0.0616"
$ws.Range("E34").Value = "  +1.01%  "

Set-TextValue "D35" "1.98"
$ws.Range("E35").Value = "  +7.26%  "

$ws.Range("E36").Value = "  +1.97%  "

$ws.Range("E37").Value = "  +13.89%  "

$ws.Range("E38").Value = "  +6.09%  "

$ws.Range("E39").Value = "  -0.14%  "

Set-TextValue "D40" "98.47"
$ws.Range("E40").Value = "  +4.08%  "

$ws.Range("E41").Value = "  +1.11%  "

Set-TextValue "D42" "1.485.05"
$ws.Range("E42").Value = "  +0.78%  "

$ws.Range("E43").Value = "  +2.96%  "

$ws.Range("E44").Value = "  +1.94%  "

Set-TextValue "D45" "2.86"
$ws.Range("E45").Value = "  +4.15%  "

$ws.Range("E46").Value = "  +0.70%  "

Set-TextValue "D47" "4.10"
$ws.Range("E47").Value = "  +19.22%  "

$ws.Range("E48").Value = "  +1.49%  "

Set-TextValue "D49" "2.98"
$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("E50").Value = "  -0.90%  "

Set-TextValue "D51" "2.251.78"
$ws.Range("E51").Value = "  +2.65%  "
